$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Data dictionary tbl_user: correct spelling "passwoord" -> "password"
$ws.Range("A7").Value = "password"

# Data dictionary tbl_ (options): fix typo "automatic" -> "automatisch"
$ws.Range("E15").Value = "hier komt automatisch de datum wanneer er een nieuwe optie is aangemaakt"

# Widen the "Gegeven" column (A) so it lines up with the other label columns
$ws.Columns.Item(1).ColumnWidth = 14.43

# Leave the cursor where the author last left it when saving
$ws.Range("O20").Select()
